$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the "Lessons learned" bullet paragraph (third bullet in the
#    numbered list at the very top of the document). Use CompareTo for
#    a case-sensitive match since -eq is case-insensitive here and the
#    document also contains a "Lessons Learned" heading further down.
# ---------------------------------------------------------------------
$lessonsPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t.CompareTo("Lessons learned") -eq 0) {
        $lessonsPara = $p
    }
}
if ($lessonsPara -ne $null) {
    $lessonsPara.Range.Delete()
}

# ---------------------------------------------------------------------
# 2) Remove the old "_GoBack" bookmark from the end of the final
#    paragraph (it will be re-added further down in the new location).
# ---------------------------------------------------------------------
try {
    $oldBookmark = $d.Bookmarks("_GoBack")
    $oldBookmark.Delete()
} catch {
}

# ---------------------------------------------------------------------
# 3) Find the empty (tab-only) paragraph that comes right after the
#    second "Interesting Items" heading near the end of the document and
#    append the two new sentences about the floor layouts.
# ---------------------------------------------------------------------
$targetPara = $null
$sawHeading = $false
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($sawHeading -eq $false) {
        if ($t.CompareTo("Interesting Items") -eq 0) {
            $sawHeading = $true
        }
    } elseif ($targetPara -eq $null) {
        $targetPara = $p
    }
}

$r = $targetPara.Range
# Range that covers the paragraph's existing content but excludes the
# trailing paragraph mark, so new runs land inside this paragraph.
$body = $d.Range($r.Start, $r.End - 1)
$body.InsertAfter("An interesting ")

$body = $d.Range($r.Start, $targetPara.Range.End - 1)
$body.Collapse(0)
$body.InsertAfter("aspect of this program is that the floor layouts used are roughly based on Adam’s house.")

$body = $d.Range($r.Start, $targetPara.Range.End - 1)
$body.Collapse(0)
$body.InsertAfter(" This was chosen since it was easy to model and younger students might better relate to a home design over an office or boxy layout.")

# ---------------------------------------------------------------------
# 4) Re-add the "_GoBack" bookmark at the end of this paragraph (right
#    after the text just inserted, before the paragraph mark).
# ---------------------------------------------------------------------
$finalRange = $targetPara.Range
$bmRange = $finalRange.Duplicate
$bmRange.Start = $finalRange.End - 1
$bmRange.End = $finalRange.End - 1
$d.Bookmarks.Add("_GoBack", $bmRange)
